$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink on D2 (and its "Hyperlink" style) before
# rewriting the data rows.
$ws.Range("D2").Hyperlinks.Delete()
$ws.Range("D2").Style = "Normal"

# Update row 2 with the new opportunity data.
$ws.Range("A2").Value = "Palermo"
$ws.Range("B2").Value = "USD 125.000"
$ws.Range("C2").Value = "Oportunidad 2 ambientes"
$ws.Range("D2").Value = "https://www.zonaprop.com.ar"

# Add row 3.
$ws.Range("A3").Value = "Recoleta"
$ws.Range("B3").Value = "USD 98.000"
$ws.Range("C3").Value = "Ideal inversión"
$ws.Range("D3").Value = "https://www.zonaprop.com.ar"

# Add row 4.
$ws.Range("A4").Value = "Belgrano"
$ws.Range("B4").Value = "USD 115.000"
$ws.Range("C4").Value = "Dueño directo impecable"
$ws.Range("D4").Value = "https://www.zonaprop.com.ar"
